$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.030704498291016
$ws.Range("B1").Value = 1.294652223587036
$ws.Range("C1").Value = 1.863481879234314
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 2.051020860671997
